# eGov - Create authtoken step checkin
# Adds a new "CreateComplaint" worksheet (after the existing "BaseAPITests"
# sheet) containing a small key/value block used to POST a "createComplaint"
# request, and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab strip (sheetId=2, rId2), matching "BaseAPITests" + "CreateComplaint".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "CreateComplaint"

# Populate B2 ("value") before B1 ("key") and A1 ("createComplaint") so the
# shared-string table grows in the same order as the source edit.
$newSheet.Range("B2").Value = "value"
$newSheet.Range("B1").Value = "key"
$newSheet.Range("A1").Value = "createComplaint"

# Auto-fit column A to its content, like Excel's "AutoFit Column Width".
$newSheet.Columns.Item(1).AutoFit() | Out-Null

# Leave the cursor on A2, and this newly-added sheet becomes the active tab.
$newSheet.Range("A2").Select()
